$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while forcing text storage
# (prevents Excel from auto-coercing numeric-looking strings, like
# "568.34", into real numbers) and without leaving any numberformat
# / style change behind on the cell.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 "63.207.83"
Set-TextValue 2 5 "  +0.41%  "

Set-TextValue 3 4 "2.552.92"
Set-TextValue 3 5 "  +3.35%  "

Set-TextValue 4 5 "  +0.03%  "

Set-TextValue 5 4 "568.34"
Set-TextValue 5 5 "  +0.45%  "

Set-TextValue 6 4 "147.17"
Set-TextValue 6 5 "  +3.28%  "

Set-TextValue 7 5 "  +0.05%  "

Set-TextValue 8 4 "0.587"
Set-TextValue 8 5 "  -0.32%  "

Set-TextValue 9 4 "2.551.00"
Set-TextValue 9 5 "  +3.33%  "

Set-TextValue 10 5 "  +0.63%  "

Set-TextValue 11 4 "5.58"
Set-TextValue 11 5 "  -1.91%  "

Set-TextValue 12 5 "  +0.49%  "

Set-TextValue 13 5 "  +0.45%  "

Set-TextValue 14 4 "27.60"
Set-TextValue 14 5 "  +3.51%  "

Set-TextValue 15 4 "3.009.51"
Set-TextValue 15 5 "  +3.42%  "

Set-TextValue 16 4 "63.128.45"
Set-TextValue 16 5 "  +0.52%  "

Set-TextValue 17 5 "  +1.98%  "

Set-TextValue 18 4 "2.554.16"
Set-TextValue 18 5 "  +3.54%  "

Set-TextValue 19 4 "11.45"
Set-TextValue 19 5 "  +1.96%  "

Set-TextValue 20 4 "335.61"
Set-TextValue 20 5 "  -1.41%  "

Set-TextValue 21 5 "  +1.57%  "

Set-TextValue 22 5 "  -0.34%  "

Set-TextValue 23 5 "  -0.17%  "

Set-TextValue 24 4 "65.27"
Set-TextValue 24 5 "  -0.50%  "

Set-TextValue 25 5 "  +9.57%  "

Set-TextValue 26 5 "  -1.78%  "

Set-TextValue 27 5 "  +6.83%  "

Set-TextValue 28 4 "8.48"
Set-TextValue 28 5 "  +4.96%  "

Set-TextValue 29 5 "  +0.04%  "

Set-TextValue 30 4 "7.37"
Set-TextValue 30 5 "  +7.61%  "

Set-TextValue 31 4 "0.0₃0822"
Set-TextValue 31 5 "  +2.60%  "

Set-TextValue 32 5 "  +0.71%  "

Set-TextValue 33 4 "176.43"
Set-TextValue 33 5 "  +0.21%  "

Set-TextValue 34 5 "  +3.85%  "

Set-TextValue 35 4 "412.68"
Set-TextValue 35 5 "  +11.94%  "

Set-TextValue 36 5 "  +1.29%  "

Set-TextValue 37 4 "18.98"
Set-TextValue 37 5 "  +1.02%  "

Set-TextValue 38 5 "  +0.53%  "

Set-TextValue 39 5 "  -0.02%  "

Set-TextValue 40 5 "  +4.28%  "

Set-TextValue 41 5 "  +0.05%  "

Set-TextValue 42 4 "39.35"
Set-TextValue 42 5 "  -3.06%  "

Set-TextValue 43 4 "153.37"
Set-TextValue 43 5 "  +2.47%  "

Set-TextValue 44 5 "  +2.45%  "

Set-TextValue 45 4 "21.15"
Set-TextValue 45 5 "  +3.21%  "

Set-TextValue 46 5 "  +0.98%  "

Set-TextValue 47 5 "  +0.56%  "

Set-TextValue 48 2 "Hedera"
Set-TextValue 48 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 48 4 "0.0524"
Set-TextValue 48 5 "  +1.77%  "

Set-TextValue 49 2 "VeChain"
Set-TextValue 49 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 49 4 "0.0239"
Set-TextValue 49 5 "  +5.91%  "

Set-TextValue 50 4 "18.41"
Set-TextValue 50 5 "  +2.53%  "

Set-TextValue 51 4 "1.78"
Set-TextValue 51 5 "  +1.66%  "

Write-Output "Applied 81 cell updates"
